$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header C1: "Gender" -> "Roll"
$ws.Range("C1").Value = "Roll"

# Row 2 data updates
$ws.Range("B2").Value = "sorna"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "kaji"
$ws.Range("E2").Value = "anam"

# Update selection to E2 (matches the diff's sheetView selection change)
$ws.Range("E2").Select()
